$d = $word.ActiveDocument

# The Title, Author and Abstract paragraphs each had their text
# fragmented across several single-word/space runs (one <w:r> per word
# and one per separating space). Word's Find/Replace operates on the
# logical paragraph text across run boundaries and, when it performs a
# replacement, collapses the matched range down to a single run -
# giving each of these three paragraphs one run carrying its full text.

$wdReplaceAll = 2
$wdFindContinue = 1

$d.Content.Find.Execute(
    "Questions: Trigonometric identities (radians)", $false, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "Questions: Trigonometric identities (radians)", $wdReplaceAll) | Out-Null

$d.Content.Find.Execute(
    "Dzhemma Ruseva", $false, $false, $false,
    $false, $false, $true, $wdFindContinue, $false,
    "Dzhemma Ruseva", $wdReplaceAll) | Out-Null

$d.Content.Find.Execute(
    "A selection of questions on trigonometric identities, where angles are measured in radians.",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "A selection of questions on trigonometric identities, where angles are measured in radians.",
    $wdReplaceAll) | Out-Null
